$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.988.72'
$ws.Range('E2').Value = '  -5.21%  '
$ws.Range('D3').Value = '2.216.11'
$ws.Range('E3').Value = '  -7.51%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'295.94"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.94%  '
$ws.Range('D6').Value = "'79.83"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -10.82%  '
$ws.Range('D7').Value = "'0.503"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.70%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.456"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.50%  '
$ws.Range('D10').Value = "'0.0772"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.11%  '
$ws.Range('D11').Value = "'27.96"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -11.23%  '
$ws.Range('D12').Value = "'45.66"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -14.08%  '
$ws.Range('D13').Value = "'0.107"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('D14').Value = '2.561.59'
$ws.Range('E14').Value = '  -7.37%  '
$ws.Range('D15').Value = "'6.08"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.46%  '
$ws.Range('D16').Value = "'13.85"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.87%  '
$ws.Range('D17').Value = '2.237.81'
$ws.Range('E17').Value = '  -6.01%  '
$ws.Range('D18').Value = "'0.708"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.69%  '
$ws.Range('D19').Value = '38.896.97'
$ws.Range('E19').Value = '  -5.22%  '
$ws.Range('E20').Value = '  -6.84%  '
$ws.Range('D21').Value = "'5.69"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.56%  '
$ws.Range('D22').Value = "'64.54"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.49%  '
$ws.Range('D23').Value = "'9.75"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.98%  '
$ws.Range('D24').Value = "'224.68"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.48%  '
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').Value = "'2.37"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -11.27%  '
$ws.Range('D27').Value = "'1.74"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.31%  '
$ws.Range('D28').Value = "'22.20"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.67%  '
$ws.Range('D29').Value = "'2.16"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('D30').Value = "'8.89"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.77%  '
$ws.Range('D31').Value = "'148.59"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.88%  '
$ws.Range('D32').Value = "'31.05"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.37%  '
$ws.Range('D33').Value = "'1.00"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').Value = "'4.72"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.24%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = "'2.33"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.98%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.0684"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.56%  '
$ws.Range('E37').Value = '  -5.55%  '
$ws.Range('D38').Value = "'2.63"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.76%  '
$ws.Range('D39').Value = "'0.0949"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.88%  '
$ws.Range('D40').Value = "'14.32"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.78%  '
$ws.Range('D41').Value = "'1.58"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.00%  '
$ws.Range('E42').Value = '  -7.72%  '
$ws.Range('D43').Value = '1.900.60'
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('D44').Value = "'2.08"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -10.93%  '
$ws.Range('D45').Value = "'0.0253"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.65%  '
$ws.Range('D46').Value = "'16.12"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.97%  '
$ws.Range('D47').Value = "'8.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.13%  '
$ws.Range('D48').Value = "'2.50"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.81%  '
$ws.Range('D49').Value = '2.429.86'
$ws.Range('E49').Value = '  -6.98%  '
$ws.Range('D50').Value = "'69.68"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.51%  '
$ws.Range('D51').Value = "'86.35"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.69%  '
